$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 41 and 42: swap coin data (RenderToken <-> Stacks), update D and E
# Leading apostrophe forces the numeric-looking price to stay text, matching
# the original inlineStr cell type instead of being coerced to a number.
$ws.Range("B41").Value = "Stacks"
$ws.Range("C41").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D41").Value = "'1.79"
$ws.Range("E41").Value = "  -1.99%  "

$ws.Range("B42").Value = "RenderToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D42").Value = "'5.16"
$ws.Range("E42").Value = "  -1.74%  "

# Update Volume(1h) percentages for all other changed rows
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("E3").Value = "  -1.69%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("E5").Value = "  -0.82%  "
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -2.15%  "
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("E10").Value = "  -1.33%  "
$ws.Range("E11").Value = "  +1.58%  "
$ws.Range("E12").Value = "  +0.52%  "
$ws.Range("E13").Value = "  +0.45%  "
$ws.Range("E14").Value = "  -0.23%  "
$ws.Range("E15").Value = "  -1.87%  "
$ws.Range("E16").Value = "  -1.28%  "
$ws.Range("E17").Value = "  -1.04%  "
$ws.Range("E18").Value = "  -2.07%  "
$ws.Range("E19").Value = "  +4.83%  "
$ws.Range("E20").Value = "  +7.89%  "
$ws.Range("E21").Value = "  -1.79%  "
$ws.Range("E22").Value = "  -1.52%  "
$ws.Range("E23").Value = "  -2.83%  "
$ws.Range("E24").Value = "  +10.72%  "
$ws.Range("E25").Value = "  -0.03%  "
$ws.Range("E26").Value = "  -4.51%  "
$ws.Range("E27").Value = "  -3.02%  "
$ws.Range("E29").Value = "  +0.18%  "
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("E31").Value = "  -1.60%  "
$ws.Range("E32").Value = "  -0.54%  "
$ws.Range("E33").Value = "  -1.40%  "
$ws.Range("E34").Value = "  -1.82%  "
$ws.Range("E35").Value = "  +5.76%  "
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("E37").Value = "  -4.51%  "
$ws.Range("E38").Value = "  -0.92%  "
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("E40").Value = "  -1.59%  "
$ws.Range("E43").Value = "  -0.08%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -3.63%  "
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("E47").Value = "  -1.49%  "
$ws.Range("E48").Value = "  -0.94%  "
$ws.Range("E49").Value = "  -1.11%  "
$ws.Range("E50").Value = "  -1.31%  "
$ws.Range("E51").Value = "  -0.42%  "
